$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.390522576735691
$ws.Range("C2").Value = 0.3045758222055497
$ws.Range("D2").Value = 0.2514654466625217
$ws.Range("F2").Value = 1.636255998397132
$ws.Range("G2").Value = 0.9425659895362344
$ws.Range("H2").Value = 0.943676326694586
$ws.Range("J2").Value = 0.3495057524152685
$ws.Range("B3").Value = 1.254972915939959
$ws.Range("C3").Value = 0.2673563839044846
$ws.Range("D3").Value = 0.2445932683339294
$ws.Range("F3").Value = 1.62635431916712
$ws.Range("G3").Value = 0.9349777015092116
$ws.Range("H3").Value = 0.9474447886092605
$ws.Range("J3").Value = 0.3382548230859754
$ws.Range("B4").Value = 1.171989245069597
$ws.Range("C4").Value = 0.2444761858400284
$ws.Range("D4").Value = 0.2404478977847759
$ws.Range("F4").Value = 1.621642553731348
$ws.Range("G4").Value = 0.9313140774567472
$ws.Range("H4").Value = 0.9504766402877465
$ws.Range("J4").Value = 0.3315679091272159
$ws.Range("B5").Value = 1.138234574412536
$ws.Range("C5").Value = 0.23514562567496
$ws.Range("D5").Value = 0.2387773606216541
$ws.Range("F5").Value = 1.620064775035715
$ws.Range("G5").Value = 0.9300697468332544
$ws.Range("H5").Value = 0.9518920319789999
$ws.Range("J5").Value = 0.3288982779646972
$ws.Range("B6").Value = 1.132633385384452
$ws.Range("C6").Value = 0.2335958943965579
$ws.Range("D6").Value = 0.2385011032533129
$ws.Range("F6").Value = 1.619823411741933
$ws.Range("G6").Value = 0.9298780925090568
$ws.Range("H6").Value = 0.9521379021899037
$ws.Range("J6").Value = 0.3284583220005857
$ws.Range("B7").Value = 1.171533766789992
$ws.Range("C7").Value = 0.2443503774998703
$ws.Range("D7").Value = 0.2404252923652592
$ws.Range("F7").Value = 1.621619891499478
$ws.Range("G7").Value = 0.9312962915995655
$ws.Range("H7").Value = 0.9504950012828317
$ws.Range("J7").Value = 0.3315316818444956
$ws.Range("B8").Value = 1.343734519601924
$ws.Range("C8").Value = 0.2917483181229557
$ws.Range("D8").Value = 0.2490805711411923
$ws.Range("F8").Value = 1.632556977901444
$ws.Range("G8").Value = 0.9397418611036272
$ws.Range("H8").Value = 0.944826271961972
$ws.Range("J8").Value = 0.3455803462165363
$ws.Range("B9").Value = 1.683360796390389
$ws.Range("C9").Value = 0.3844763767526729
$ws.Range("D9").Value = 0.266639606090294
$ws.Range("F9").Value = 1.6649372110992
$ws.Range("G9").Value = 0.9642836860957402
$ws.Range("H9").Value = 0.9394376334646779
$ws.Range("J9").Value = 0.3748995633946919
$ws.Range("B10").Value = 1.934099774438209
$ws.Range("C10").Value = 0.4524745006963826
$ws.Range("D10").Value = 0.2798958999282632
$ws.Range("F10").Value = 1.695506466407551
$ws.Range("G10").Value = 0.9872939869356401
$ws.Range("H10").Value = 0.9390142786266154
$ws.Range("J10").Value = 0.3975424788065141
$ws.Range("B11").Value = 2.048440787153311
$ws.Range("C11").Value = 0.4833821916792544
$ws.Range("D11").Value = 0.2860035794087992
$ws.Range("F11").Value = 1.710910483712809
$ws.Range("G11").Value = 0.9988681519160139
$ws.Range("H11").Value = 0.9395990306898057
$ws.Range("J11").Value = 0.4080878588138006
$ws.Range("B12").Value = 2.091779057084409
$ws.Range("C12").Value = 0.4950825310036748
$ws.Range("D12").Value = 0.2883274727263085
$ws.Range("F12").Value = 1.716960938625007
$ws.Range("G12").Value = 1.003412092739637
$ws.Range("H12").Value = 0.9399330080993877
$ws.Range("J12").Value = 0.4121167255259621
$ws.Range("B13").Value = 2.082443618181401
$ws.Range("C13").Value = 0.4925628216400924
$ws.Range("D13").Value = 0.2878264904803416
$ws.Range("F13").Value = 1.715648170728571
$ws.Range("G13").Value = 1.002426280131743
$ws.Range("H13").Value = 0.9398560624146057
$ws.Range("J13").Value = 0.4112474501547041
$ws.Range("B14").Value = 2.052005456127631
$ws.Range("C14").Value = 0.4843448611820804
$ws.Range("D14").Value = 0.2861945464987343
$ws.Range("F14").Value = 1.711403892635872
$ws.Range("G14").Value = 0.9992387458431153
$ws.Range("H14").Value = 0.9396242472583936
$ws.Range("J14").Value = 0.4084186015218165
$ws.Range("B15").Value = 2.033366380747566
$ws.Range("C15").Value = 0.4793106331741797
$ws.Range("D15").Value = 0.2851963703479186
$ws.Range("F15").Value = 1.708832503361066
$ws.Range("G15").Value = 0.9973073203651666
$ws.Range("H15").Value = 0.9394969331623884
$ws.Range("J15").Value = 0.406690490885353
$ws.Range("B16").Value = 1.92663290380807
$ws.Range("C16").Value = 0.4504540953020069
$ws.Range("D16").Value = 0.2794982984949854
$ws.Range("F16").Value = 1.694530080725514
$ws.Range("G16").Value = 0.9865600207638465
$ws.Range("H16").Value = 0.9389917709470978
$ws.Range("J16").Value = 0.3968582715288278
$ws.Range("B17").Value = 1.861226533663967
$ws.Range("C17").Value = 0.4327450392742662
$ws.Range("D17").Value = 0.276022468822589
$ws.Range("F17").Value = 1.68614099935364
$ws.Range("G17").Value = 0.9802517061977767
$ws.Range("H17").Value = 0.9388815074020158
$ws.Range("J17").Value = 0.3908895072477208
$ws.Range("B18").Value = 1.823632754164237
$ws.Range("C18").Value = 0.4225569132182727
$ws.Range("D18").Value = 0.2740305457979844
$ws.Range("F18").Value = 1.681456698275412
$ws.Range("G18").Value = 0.9767273797274925
$ws.Range("H18").Value = 0.938891207520868
$ws.Range("J18").Value = 0.3874794728863122
$ws.Range("B19").Value = 1.810908660638233
$ws.Range("C19").Value = 0.4191069903943685
$ws.Range("D19").Value = 0.2733573678571588
$ws.Range("F19").Value = 1.679894812705669
$ws.Range("G19").Value = 0.9755519155593504
$ws.Range("H19").Value = 0.9389070271675166
$ws.Range("J19").Value = 0.3863288412186847
$ws.Range("B20").Value = 1.868186440194961
$ws.Range("C20").Value = 0.4346304422927574
$ws.Range("D20").Value = 0.2763917236641049
$ws.Range("F20").Value = 1.68701943757479
$ws.Range("G20").Value = 0.9809124545064947
$ws.Range("H20").Value = 0.9388856719295688
$ws.Range("J20").Value = 0.3915225053699629
$ws.Range("B21").Value = 2.060944803496852
$ws.Range("C21").Value = 0.4867587777061431
$ws.Range("D21").Value = 0.2866735887297125
$ws.Range("F21").Value = 1.712644627736196
$ws.Range("G21").Value = 1.00017061558674
$ws.Range("H21").Value = 0.9396892764472113
$ws.Range("J21").Value = 0.409248535198742
$ws.Range("B22").Value = 2.187155496186222
$ws.Range("C22").Value = 0.5208058482789966
$ws.Range("D22").Value = 0.2934577451017901
$ws.Range("F22").Value = 1.730659577181953
$ws.Range("G22").Value = 1.013696590243853
$ws.Range("H22").Value = 0.9408708322607708
$ws.Range("J22").Value = 0.4210409486082796
$ws.Range("B23").Value = 2.119773320903505
$ws.Range("C23").Value = 0.5026363236398197
$ws.Range("D23").Value = 0.2898310467702743
$ws.Range("F23").Value = 1.720928063963385
$ws.Range("G23").Value = 1.006390908299068
$ws.Range("H23").Value = 0.9401799034073406
$ws.Range("J23").Value = 0.4147280259293211
$ws.Range("B24").Value = 1.865039837736902
$ws.Range("C24").Value = 0.4337780732437295
$ws.Range("D24").Value = 0.2762247636809008
$ws.Range("F24").Value = 1.686621863779465
$ws.Range("G24").Value = 0.9806134113526355
$ws.Range("H24").Value = 0.9388835615247331
$ws.Range("J24").Value = 0.3912362600170383
$ws.Range("B25").Value = 1.591271924357443
$ws.Range("C25").Value = 0.3594140354481397
$ws.Range("D25").Value = 0.261826857838571
$ws.Range("F25").Value = 1.654995272167582
$ws.Range("G25").Value = 0.9567783352778321
$ws.Range("H25").Value = 0.9402776886121273
$ws.Range("J25").Value = 0.3382548230859754
